$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on each touched cell individually (multi-area NumberFormat only applies to first area)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("B48").NumberFormat = "@"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("B49").NumberFormat = "@"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("B50").NumberFormat = "@"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("E51").NumberFormat = "@"

# Apply new values
$ws.Range("D2").Value = '45.584.04'
$ws.Range("E2").Value = '  -2.23%  '
$ws.Range("D3").Value = '2.367.56'
$ws.Range("E3").Value = '  +3.14%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '298.39'
$ws.Range("E5").Value = '  -1.96%  '
$ws.Range("D6").Value = '96.92'
$ws.Range("E6").Value = '  -3.94%  '
$ws.Range("D7").Value = '0.560'
$ws.Range("E7").Value = '  -1.27%  '
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("D9").Value = '0.500'
$ws.Range("E9").Value = '  -4.17%  '
$ws.Range("D10").Value = '33.83'
$ws.Range("E10").Value = '  -8.01%  '
$ws.Range("D11").Value = '0.0782'
$ws.Range("E11").Value = '  -1.33%  '
$ws.Range("D12").Value = '7.00'
$ws.Range("E12").Value = '  -5.77%  '
$ws.Range("E13").Value = '  +0.30%  '
$ws.Range("D14").Value = '2.729.18'
$ws.Range("E14").Value = '  +3.13%  '
$ws.Range("D15").Value = '2.369.41'
$ws.Range("E15").Value = '  +2.87%  '
$ws.Range("D16").Value = '0.813'
$ws.Range("E16").Value = '  -0.16%  '
$ws.Range("E17").Value = '  -1.93%  '
$ws.Range("D18").Value = '45.553.09'
$ws.Range("E18").Value = '  -2.30%  '
$ws.Range("D19").Value = '12.63'
$ws.Range("E19").Value = '  -4.62%  '
$ws.Range("D20").Value = '0.0₃0937'
$ws.Range("E20").Value = '  -0.79%  '
$ws.Range("E21").Value = '  -1.29%  '
$ws.Range("D22").Value = '66.59'
$ws.Range("E22").Value = '  -0.08%  '
$ws.Range("D23").Value = '241.14'
$ws.Range("E23").Value = '  -3.28%  '
$ws.Range("D24").Value = '2.74'
$ws.Range("E24").Value = '  -6.06%  '
$ws.Range("D25").Value = '0.999'
$ws.Range("E25").Value = '  +0.03%  '
$ws.Range("E26").Value = '  -2.72%  '
$ws.Range("D27").Value = '38.27'
$ws.Range("E27").Value = '  -11.02%  '
$ws.Range("D28").Value = '2.20'
$ws.Range("E28").Value = '  -3.08%  '
$ws.Range("D29").Value = '9.59'
$ws.Range("E29").Value = '  -2.96%  '
$ws.Range("D30").Value = '3.76'
$ws.Range("E30").Value = '  +16.46%  '
$ws.Range("D31").Value = '20.77'
$ws.Range("E31").Value = '  +3.37%  '
$ws.Range("D32").Value = '2.74'
$ws.Range("E32").Value = '  -2.22%  '
$ws.Range("D33").Value = '146.41'
$ws.Range("E33").Value = '  -0.57%  '
$ws.Range("D34").Value = '5.45'
$ws.Range("E34").Value = '  -4.26%  '
$ws.Range("E35").Value = '  -4.40%  '
$ws.Range("D36").Value = '0.111'
$ws.Range("E36").Value = '  -2.86%  '
$ws.Range("D37").Value = '1.89'
$ws.Range("E37").Value = '  +6.06%  '
$ws.Range("D38").Value = '0.115'
$ws.Range("E38").Value = '  -2.70%  '
$ws.Range("D39").Value = '15.15'
$ws.Range("E39").Value = '  -5.59%  '
$ws.Range("D40").Value = '3.80'
$ws.Range("E40").Value = '  -7.13%  '
$ws.Range("E41").Value = '  -2.99%  '
$ws.Range("D42").Value = '3.18'
$ws.Range("E42").Value = '  -6.90%  '
$ws.Range("D43").Value = '1.938.16'
$ws.Range("E43").Value = '  +4.73%  '
$ws.Range("D44").Value = '1.00'
$ws.Range("E44").Value = '  -0.01%  '
$ws.Range("D45").Value = '92.38'
$ws.Range("E45").Value = '  +5.60%  '
$ws.Range("D46").Value = '1.76'
$ws.Range("E46").Value = '  -11.13%  '
$ws.Range("D47").Value = '8.46'
$ws.Range("E47").Value = '  +7.07%  '
$ws.Range("B48").Value = 'RocketPoolETH'
$ws.Range("C48").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D48").Value = '2.599.18'
$ws.Range("E48").Value = '  +3.09%  '
$ws.Range("B49").Value = 'Aave'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D49").Value = '97.80'
$ws.Range("E49").Value = '  +1.65%  '
$ws.Range("B50").Value = 'Algorand'
$ws.Range("C50").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D50").Value = '0.182'
$ws.Range("E50").Value = '  -6.71%  '
$ws.Range("D51").Value = '68.18'
$ws.Range("E51").Value = '  -7.47%  '
